# Auto-generated edit script: applies Aegis_Profits market-data refresh
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101498

$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -307488

$ws.Range("H92").Value = 878.625
$ws.Range("I92").Value = 730.4375
$ws.Range("J92").Value = 1175
$ws.Range("K92").Value = 730.4375
$ws.Range("L92").Value = 1175
$ws.Range("M92").Value = 517.5625
$ws.Range("N92").Value = -3671

$ws.Range("H97").Value = 1300
$ws.Range("J97").Value = 1300
$ws.Range("L97").Value = 3900
$ws.Range("N97").Value = -4892

$ws.Range("H98").Value = 1131.5714
$ws.Range("I98").Value = 1254.5
$ws.Range("J98").Value = 394
$ws.Range("K98").Value = 1254.5
$ws.Range("L98").Value = 394
$ws.Range("M98").Value = 243.5
$ws.Range("N98").Value = -3390

$ws.Range("H99").Value = 11407.777
$ws.Range("I99").Value = 12740
$ws.Range("J99").Value = 750
$ws.Range("K99").Value = 38220
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -36722
$ws.Range("N99").Value = -5246

$ws.Range("H116").Value = 3196.6667
$ws.Range("I116").Value = 2633.3333
$ws.Range("J116").Value = 3760
$ws.Range("K116").Value = 2633.3333
$ws.Range("L116").Value = 3760
$ws.Range("M116").Value = 808.6667000000002
$ws.Range("N116").Value = -10644

$ws.Range("H122").Value = 1131.5714
$ws.Range("I122").Value = 1254.5
$ws.Range("J122").Value = 394
$ws.Range("K122").Value = 3763.5
$ws.Range("L122").Value = 1182
$ws.Range("M122").Value = -1313.5
$ws.Range("N122").Value = -6082

$ws.Range("H132").Value = 8276.23
$ws.Range("I132").Value = 9225.348
$ws.Range("J132").Value = 999.6667
$ws.Range("K132").Value = 27676.044
$ws.Range("L132").Value = 2999.0001
$ws.Range("M132").Value = -25146.044
$ws.Range("N132").Value = -8059.0001

$ws.Range("H135").Value = 1115.3864
$ws.Range("I135").Value = 631.4762
$ws.Range("J135").Value = 1557.2174
$ws.Range("K135").Value = 5683.2858
$ws.Range("L135").Value = 14014.9566
$ws.Range("M135").Value = -3148.2858
$ws.Range("N135").Value = -19084.9566

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6250
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885

$ws.Range("H44").Value = 11348.833
$ws.Range("J44").Value = 13009.8
$ws.Range("L44").Value = 13009.8
$ws.Range("N44").Value = -13985.8

$ws.Range("H47").Value = 9500
$ws.Range("J47").Value = 9500
$ws.Range("L47").Value = 9500
$ws.Range("N47").Value = -10950

$ws.Range("H55").Value = 12008.833
$ws.Range("J55").Value = 13010.6
$ws.Range("L55").Value = 13010.6
$ws.Range("N55").Value = -13640.6

$ws.Range("H61").Value = 2247.0667
$ws.Range("I61").Value = 2122.4443
$ws.Range("J61").Value = 2434
$ws.Range("K61").Value = 2122.4443
$ws.Range("L61").Value = 2434
$ws.Range("M61").Value = -1910.4443
$ws.Range("N61").Value = -2858

$ws.Range("H132").Value = 3061.9736
$ws.Range("I132").Value = 2900.0667
$ws.Range("J132").Value = 3669.125
$ws.Range("K132").Value = 8700.2001
$ws.Range("L132").Value = 11007.375
$ws.Range("M132").Value = -6170.2001
$ws.Range("N132").Value = -16067.375

$ws.Range("H136").Value = 2247.0667
$ws.Range("I136").Value = 2122.4443
$ws.Range("J136").Value = 2434
$ws.Range("K136").Value = 6367.3329
$ws.Range("L136").Value = 7302
$ws.Range("M136").Value = -3817.3329
$ws.Range("N136").Value = -12402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 33333.332
$ws.Range("J114").Value = 33333.332
$ws.Range("L114").Value = 33333.332
$ws.Range("N114").Value = -42011.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36132

$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20952

$ws.Range("H134").Value = 1795.4166
$ws.Range("I134").Value = 1555.625
$ws.Range("J134").Value = 2275
$ws.Range("K134").Value = 4666.875
$ws.Range("L134").Value = 6825
$ws.Range("M134").Value = -2131.875
$ws.Range("N134").Value = -11895

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1901
$ws.Range("I68").Value = 400
$ws.Range("J68").Value = 2651.5
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 7954.5
$ws.Range("M68").Value = -389
$ws.Range("N68").Value = -9576.5

$ws.Range("H71").Value = 1901
$ws.Range("I71").Value = 400
$ws.Range("J71").Value = 2651.5
$ws.Range("K71").Value = 3600
$ws.Range("L71").Value = 23863.5
$ws.Range("M71").Value = 456
$ws.Range("N71").Value = -31975.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 100000
$ws.Range("I20").Value = 100000
$ws.Range("K20").Value = 100000
$ws.Range("M20").Value = -99755

$ws.Range("H80").Value = 55614892
$ws.Range("I80").Value = 100104620
$ws.Range("J80").Value = 2726.25
$ws.Range("K80").Value = 100104620
$ws.Range("L80").Value = 2726.25
$ws.Range("M80").Value = -100103622
$ws.Range("N80").Value = -4722.25

$ws.Range("H83").Value = 55614892
$ws.Range("I83").Value = 100104620
$ws.Range("J83").Value = 2726.25
$ws.Range("K83").Value = 500523100
$ws.Range("L83").Value = 13631.25
$ws.Range("M83").Value = -500518108
$ws.Range("N83").Value = -23615.25

$ws.Range("H132").Value = 3266.0557
$ws.Range("I132").Value = 3088.6155
$ws.Range("J132").Value = 3727.4
$ws.Range("K132").Value = 9265.8465
$ws.Range("L132").Value = 11182.2
$ws.Range("M132").Value = -6735.8465
$ws.Range("N132").Value = -16242.2

$ws.Range("H141").Value = 47495
$ws.Range("J141").Value = 47495
$ws.Range("L141").Value = 47495
$ws.Range("N141").Value = -57855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2139.9
$ws.Range("I136").Value = 1922.1111
$ws.Range("J136").Value = 4100
$ws.Range("K136").Value = 5766.3333
$ws.Range("L136").Value = 12300
$ws.Range("M136").Value = -3216.3333
$ws.Range("N136").Value = -17400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2331.111
$ws.Range("I126").Value = 2245
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 6735
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -4265
$ws.Range("N126").Value = -12140

$ws.Range("H140").Value = 62921.43
$ws.Range("J140").Value = 62921.43
$ws.Range("L140").Value = 62921.43
$ws.Range("N140").Value = -73281.43

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

